$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.736.44"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.640.80"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "217.58"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "19.12"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.870.37"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "1.639.62"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "64.68"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "26.733.52"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "214.78"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("E22").Value = "  +6.56%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "145.55"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "1.287.42"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "0.816"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").Value = "1.780.54"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "61.29"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "7.66"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("E51").Value = "  +0.08%  "